$wb = $excel.ActiveWorkbook

# --- Update Generic sheet: NrBuckets 4 -> 5 ---
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B4").Value = 5

# --- Add bucket row 6 to ForecastedAverageDemand sheet ---
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvg.Range("A2:N2").Copy($wsAvg.Range("A6:N6"))
$wsAvg.Range("A6").Value = 4
$wsAvg.Range("B6:F6").Value = 0
$wsAvg.Range("G6").Value = 4000
$wsAvg.Range("H6").Value = 2000
$wsAvg.Range("I6").Value = 8000
$wsAvg.Range("J6").Value = 700
$wsAvg.Range("K6:N6").Value = 0

# --- Add bucket row 6 to ForcastedStandardDeviation sheet ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStd.Range("A2:N2").Copy($wsStd.Range("A6:N6"))
$wsStd.Range("A6").Value = 4
$wsStd.Range("B6:F6").Value = 0
$wsStd.Range("G6").Value = 7500
$wsStd.Range("H6").Value = 3500
$wsStd.Range("I6").Value = 9000
$wsStd.Range("J6").Value = 1400
$wsStd.Range("K6:N6").Value = 0
